# Updates the "cryptos" worksheet with refreshed price/volume figures
# (and a RenderToken / RocketPoolETH row swap), matching a scheduled
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of (cell, newValue) pairs. Values are plain strings - many of them
# look numeric/date-like (e.g. "241.77", "1.10", "36.448.16") but must be
# stored as literal text, exactly as the source data feed produced them.
$updates = @(
    @('D2', '36.448.16'),
    @('E2', '  -2.02%  '),
    @('D3', '2.054.05'),
    @('E3', '  -0.57%  '),
    @('E4', '  -0.08%  '),
    @('D5', '241.77'),
    @('E5', '  -3.02%  '),
    @('D6', '0.664'),
    @('E6', '  -0.49%  '),
    @('E7', '  +0.01%  '),
    @('D8', '54.33'),
    @('E8', '  -7.17%  '),
    @('D9', '58.32'),
    @('E9', '  -3.70%  '),
    @('E10', '  -7.59%  '),
    @('E11', '  -5.40%  '),
    @('D12', '0.106'),
    @('E12', '  -3.08%  '),
    @('D13', '0.891'),
    @('E13', '  -2.80%  '),
    @('D14', '14.58'),
    @('E14', '  -8.14%  '),
    @('E15', '  -0.36%  '),
    @('D16', '5.33'),
    @('E16', '  -9.18%  '),
    @('D17', '2.055.74'),
    @('E17', '  -0.64%  '),
    @('D18', '36.395.91'),
    @('E18', '  -2.20%  '),
    @('D19', '16.64'),
    @('E19', '  -10.38%  '),
    @('D20', '71.87'),
    @('E21', '  -6.37%  '),
    @('D22', '237.66'),
    @('E22', '  -0.65%  '),
    @('D23', '5.23'),
    @('E23', '  -5.54%  '),
    @('E24', '  +0.21%  '),
    @('D25', '2.34'),
    @('E25', '  -5.77%  '),
    @('D26', '9.33'),
    @('E26', '  -3.56%  '),
    @('D27', '2.11'),
    @('E27', '  -5.34%  '),
    @('D28', '162.72'),
    @('E28', '  -5.10%  '),
    @('E29', '  -1.24%  '),
    @('E30', '  -3.32%  '),
    @('D31', '5.02'),
    @('E31', '  -10.56%  '),
    @('E32', '  -0.18%  '),
    @('E33', '  -8.70%  '),
    @('E34', '  -6.28%  '),
    @('E35', '  +0.01%  '),
    @('E36', '  +1.23%  '),
    @('E37', '  -6.46%  '),
    @('E38', '  -7.48%  '),
    @('D39', '1.23'),
    @('E39', '  -8.65%  '),
    @('E41', '  -6.00%  '),
    @('E42', '  -9.53%  '),
    @('D43', '1.10'),
    @('E43', '  -5.87%  '),
    @('D44', '93.25'),
    @('E44', '  -7.37%  '),
    @('D45', '0.0895'),
    @('E45', '  -12.02%  '),
    @('D46', '1.379.75'),
    @('E46', '  +4.86%  '),
    @('D47', '15.61'),
    @('E47', '  -11.24%  '),
    @('E48', '  +3.68%  '),
    @('E49', '  -1.12%  '),
    @('B50', 'RocketPoolETH'),
    @('C50', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'),
    @('D50', '2.245.10'),
    @('E50', '  -0.39%  '),
    @('B51', 'RenderToken'),
    @('C51', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @('D51', '2.25'),
    @('E51', '  -8.02%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $cell = $ws.Range($cellRef)
    # Prefix with an apostrophe so the COM layer always treats the value as
    # literal text instead of auto-converting number/date-looking strings.
    $cell.Value = "'" + $newVal
    # Re-apply the default "Normal" style so no stray number-format / style
    # index is left behind on the cell (keeps cell formatting unchanged).
    $cell.Style = "Normal"
}

# D21 contains a subscript-three Unicode character (U+2083) inside the
# price text. Build it with the format operator so the scripting engine
# doesn't try to numerically evaluate the pieces being combined.
$d21Value = "0.0{0}0853" -f [char]0x2083
$d21 = $ws.Range("D21")
$d21.Value = "'" + $d21Value
$d21.Style = "Normal"
